$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        if ($t.Contains($needle)) {
            return $i
        }
    }
    return -1
}

function Assert-Found($ok, [string]$label) {
    if (-not $ok) {
        throw "Find.Execute failed to match: $label"
    }
}

# ---------------------------------------------------------------------------
# 1. "we'll a statistic" -> "we'll work with a statistic"
#    The whole sentence is a single run, so scope the Find to exactly that
#    run's span (start of its paragraph, length of the old text) to avoid
#    bleeding the edit into neighboring runs.
# ---------------------------------------------------------------------------
$needle1 = ("we" + [char]8217 + "ll a statistic called R")
$idx1 = Get-ParagraphIndexContaining $d $needle1
if ($idx1 -lt 0) { throw "Could not locate paragraph for edit #1" }
$old1 = ("That" + [char]8217 + "s what this lesson is about. As you" + [char]8217 + "ll see, we" + [char]8217 + "ll a statistic called R with the long-winded name")
$new1 = ("That" + [char]8217 + "s what this lesson is about. As you" + [char]8217 + "ll see, we" + [char]8217 + "ll work with a statistic called R with the long-winded name")
$p1 = $d.Paragraphs($idx1).Range
$rng1 = $d.Range($p1.Start, $p1.Start + $old1.Length)
$ok1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Assert-Found $ok1 "edit #1 (we'll a statistic -> we'll work with a statistic)"

# ---------------------------------------------------------------------------
# 2. Bullet about standard deviation:
#      "The | standard deviation | us another, usually better, way ..."
#    becomes
#      "Another way to quantify the spread of the red bars is with the |
#       standard deviation | . The black I-shaped mark ..."
#    Do this as two narrowly-scoped edits so the italic "standard deviation"
#    run is left completely untouched.
# ---------------------------------------------------------------------------
$idx2 = Get-ParagraphIndexContaining $d "usually better, way to measure"
if ($idx2 -lt 0) { throw "Could not locate paragraph for edit #2" }

# 2a. First run of the paragraph: "The" -> "Another way ... with the"
$p2 = $d.Paragraphs($idx2).Range
$rng2a = $d.Range($p2.Start, $p2.Start + "The".Length)
$ok2a = $rng2a.Find.Execute("The", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Another way to quantify the spread of the red bars is with the", 2)
Assert-Found $ok2a "edit #2a (The -> Another way to quantify...)"

# 2b. Drop the redundant " " run and rewrite the tail run that follows the
#     italic "standard deviation" run.
$p2b = $d.Paragraphs($idx2).Range
$old2b = " us another, usually better, way to measure the spread of the red bars is with . The black I-shaped mark spans a vertical distance of one standard deviation."
$new2b = ". The black I-shaped mark spans a vertical distance of one standard deviation."
$ok2b = $p2b.Find.Execute($old2b, $true, $false, $false, $false, $false, $true, 1, $false, $new2b, 2)
Assert-Found $ok2b "edit #2b (drop redundant run / rewrite tail)"

# ---------------------------------------------------------------------------
# 3. typo fix: "reponse variable" -> "response variable"
#    This run runs to the end of its paragraph, so a whole-paragraph-scoped
#    Find is safe (no following run to accidentally coalesce into).
# ---------------------------------------------------------------------------
$idx3 = Get-ParagraphIndexContaining $d "partial explanation for the reponse variable"
if ($idx3 -lt 0) { throw "Could not locate paragraph for edit #3" }
$p3 = $d.Paragraphs($idx3).Range
$ok3 = $p3.Find.Execute("partial explanation for the reponse variable.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "partial explanation for the response variable.", 2)
Assert-Found $ok3 "edit #3 (reponse -> response)"

# ---------------------------------------------------------------------------
# 4. typo fix: "expanatory variable" -> "explanatory variable"
# ---------------------------------------------------------------------------
$idx4 = Get-ParagraphIndexContaining $d "translate the expanatory variable"
if ($idx4 -lt 0) { throw "Could not locate paragraph for edit #4" }
$p4 = $d.Paragraphs($idx4).Range
$ok4 = $p4.Find.Execute("translate the expanatory variable into a value", $true, $false, $false, $false, $false, $true, 1, $false, `
    "translate the explanatory variable into a value", 2)
Assert-Found $ok4 "edit #4 (expanatory -> explanatory)"

Write-Output "All 4 edits applied successfully."
